$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Simple value updates (price / volume columns)
$ws.Range("D2").Value = '67.651.18'
$ws.Range("E2").Value = '  +0.85%  '
$ws.Range("D3").Value = '3.497.73'
$ws.Range("E3").Value = '  -0.17%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '597.36'
$ws.Range("E5").Value = '  +0.41%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '179.62'
$ws.Range("E6").Value = '  +3.83%  '
$ws.Range("E7").Value = '  +0.11%  '
$ws.Range("D8").Value = '3.496.43'
$ws.Range("E8").Value = '  -0.18%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.592'
$ws.Range("E9").Value = '  -1.40%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.139'
$ws.Range("E10").Value = '  +6.29%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.12'
$ws.Range("E11").Value = '  -2.24%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.432'
$ws.Range("E12").Value = '  -0.22%  '
$ws.Range("D13").Value = '4.104.51'
$ws.Range("E13").Value = '  -0.12%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '32.23'
$ws.Range("E14").Value = '  +11.11%  '
$ws.Range("E15").Value = '  +1.76%  '
$ws.Range("D16").Value = '67.654.39'
$ws.Range("E16").Value = '  +0.86%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000177'
$ws.Range("E17").Value = '  -0.64%  '
$ws.Range("D18").Value = '3.504.62'
$ws.Range("E18").Value = '  -0.52%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.29'
$ws.Range("E19").Value = '  -0.32%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.36'
$ws.Range("E20").Value = '  +1.52%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '395.41'
$ws.Range("E21").Value = '  +0.28%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.89'
$ws.Range("E22").Value = '  -1.54%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '73.20'
$ws.Range("E23").Value = '  +0.10%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.537'
$ws.Range("E25").Value = '  +0.17%  '
$ws.Range("E26").Value = '  +0.35%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0000121'
$ws.Range("E27").Value = '  +0.27%  '
$ws.Range("E28").Value = '  +2.32%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.177'
$ws.Range("E29").Value = '  -2.50%  '
$ws.Range("E30").Value = '  +0.13%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.18'
$ws.Range("E31").Value = '  -1.71%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.43'
$ws.Range("E32").Value = '  -0.78%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.07'
$ws.Range("E33").Value = '  +0.70%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '23.67'
$ws.Range("E34").Value = '  -0.33%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '7.40'
$ws.Range("E35").Value = '  +0.55%  '
$ws.Range("E37").Value = '  -3.03%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '163.99'
$ws.Range("E38").Value = '  +0.32%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.875'
$ws.Range("E39").Value = '  -0.63%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.90'
$ws.Range("E40").Value = '  +0.16%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '4.66'
$ws.Range("E43").Value = '  -0.54%  '
$ws.Range("D44").Value = '2.853.06'
$ws.Range("E44").Value = '  +1.52%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '26.48'
$ws.Range("E47").Value = '  -3.07%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '42.06'
$ws.Range("E48").Value = '  -1.55%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0301'
$ws.Range("E49").Value = '  -0.34%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '338.33'
$ws.Range("E50").Value = '  +0.67%  '
$ws.Range("E51").Value = '  -1.84%  '

# Rows re-ranked: coin identity (name/link) and stats moved between rows
$ws.Range("B41").Value = 'RenderToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.92'
$ws.Range("E41").Value = '  -0.59%  '
$ws.Range("B42").Value = 'dogwifhat'
$ws.Range("C42").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.73'
$ws.Range("E42").Value = '  +6.42%  '
$ws.Range("B45").Value = 'Hedera'
$ws.Range("C45").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0727'
$ws.Range("E45").Value = '  -2.41%  '
$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '26.01'
$ws.Range("E46").Value = '  -1.24%  '
